$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new person type row: ID 4 = "retired non working adult"
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "retired non working adult"

# Leave the selection where the author left it before saving
$ws.Range("B7").Select()
